$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '''27.974.55'
$ws.Range("E2").Value = '  +0.71%  '
$ws.Range("D3").Value = '''1.813.68'
$ws.Range("E3").Value = '  +1.81%  '
$ws.Range("D4").Value = '''1.002'
$ws.Range("E4").Value = '  -0.18%  '
$ws.Range("D5").Value = '''310.51'
$ws.Range("E5").Value = '  -0.10%  '
$ws.Range("D6").Value = '''1.001'
$ws.Range("E6").Value = '  -0.16%  '
$ws.Range("D7").Value = '''0.4972'
$ws.Range("E7").Value = '  -2.87%  '
$ws.Range("D8").Value = '''0.3890'
$ws.Range("E8").Value = '  +3.11%  '
$ws.Range("D9").Value = '''0.09669'
$ws.Range("E9").Value = '  +24.59%  '
$ws.Range("D10").Value = '''1.103'
$ws.Range("E10").Value = '  +1.61%  '
$ws.Range("D11").Value = '''41.06'
$ws.Range("E11").Value = '  -0.24%  '
$ws.Range("D12").Value = '''6.446'
$ws.Range("E12").Value = '  +3.96%  '
$ws.Range("E13").Value = '  +1.93%  '
$ws.Range("D14").Value = '''1.002'
$ws.Range("E14").Value = '  -0.16%  '
$ws.Range("D15").Value = '''1.812.86'
$ws.Range("E15").Value = '  +2.04%  '
$ws.Range("D16").Value = '''7.308'
$ws.Range("E16").Value = '  +1.99%  '
$ws.Range("D17").Value = '''0.00001127'
$ws.Range("E17").Value = '  +5.32%  '
$ws.Range("D18").Value = '''92.64'
$ws.Range("E18").Value = '  +0.75%  '
$ws.Range("D19").Value = '''0.06637'
$ws.Range("E19").Value = '  +1.48%  '
$ws.Range("D20").Value = '''1.002'
$ws.Range("E20").Value = '  -0.12%  '
$ws.Range("D21").Value = '''17.11'
$ws.Range("E21").Value = '  +0.93%  '
$ws.Range("E22").Value = '  +0.06%  '
$ws.Range("D23").Value = '''28.017.51'
$ws.Range("E23").Value = '  +0.69%  '
$ws.Range("D24").Value = '''11.12'
$ws.Range("E24").Value = '  +1.50%  '
$ws.Range("D25").Value = '''2.246'
$ws.Range("E25").Value = '  +0.14%  '
$ws.Range("D26").Value = '''159.49'
$ws.Range("E26").Value = '  +0.57%  '
$ws.Range("B27").Value = 'EthereumClassic'
$ws.Range("C27").Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$ws.Range("D27").Value = '''20.62'
$ws.Range("E27").Value = '  +2.17%  '
$ws.Range("B28").Value = 'WrappedliquidstakedEther2.0'
$ws.Range("C28").Value = 'https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth'
$ws.Range("D28").Value = '''2.021.33'
$ws.Range("E28").Value = '  +1.87%  '
$ws.Range("D29").Value = '''2.395'
$ws.Range("E29").Value = '  +2.04%  '
$ws.Range("D30").Value = '''127.91'
$ws.Range("E30").Value = '  +2.87%  '
$ws.Range("D31").Value = '''0.1061'
$ws.Range("E31").Value = '  -1.44%  '
$ws.Range("D32").Value = '''1.040'
$ws.Range("E32").Value = '  +1.11%  '
$ws.Range("D33").Value = '''5.577'
$ws.Range("E33").Value = '  +1.79%  '
$ws.Range("E34").Value = '  +0.36%  '
$ws.Range("D35").Value = '''0.06721'
$ws.Range("E35").Value = '  -4.65%  '
$ws.Range("D36").Value = '''8.952'
$ws.Range("E36").Value = '  +2.80%  '
$ws.Range("D37").Value = '''0.02334'
$ws.Range("E37").Value = '  +1.36%  '
$ws.Range("D38").Value = '''0.2141'
$ws.Range("E38").Value = '  +1.19%  '
$ws.Range("D39").Value = '''4.946'
$ws.Range("E39").Value = '  -1.16%  '
$ws.Range("E40").Value = '  -2.22%  '
$ws.Range("D41").Value = '''0.6187'
$ws.Range("E41").Value = '  +1.64%  '
$ws.Range("E42").Value = '  -0.10%  '
$ws.Range("D43").Value = '''1.148'
$ws.Range("E43").Value = '  +0.03%  '
$ws.Range("D44").Value = '''13.10'
$ws.Range("E44").Value = '  +0.56%  '
$ws.Range("D45").Value = '''0.5875'
$ws.Range("E45").Value = '  -1.31%  '
$ws.Range("D46").Value = '''3.693'
$ws.Range("E46").Value = '  -0.60%  '
$ws.Range("D47").Value = '''1.277'
$ws.Range("E47").Value = '  -2.99%  '
$ws.Range("D48").Value = '''122.82'
$ws.Range("E48").Value = '  -3.71%  '
$ws.Range("D49").Value = '''1.935'
$ws.Range("E49").Value = '  +2.17%  '
$ws.Range("D50").Value = '''1.178'
$ws.Range("E50").Value = '  -2.52%  '
$ws.Range("D51").Value = '''0.06786'
$ws.Range("E51").Value = '  +1.23%  '
